$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add the "Normal (Web)" paragraph style (referenced by the new citation
#    paragraphs below) with the same formatting the target document uses.
# ---------------------------------------------------------------------------
$normalWeb = $d.Styles.Add("Normal (Web)")
$normalWeb.BaseStyle = $d.Styles.Item("Normal")
$normalWeb.Priority = 99
$normalWeb.UnhideWhenUsed = $true

$pf = $normalWeb.ParagraphFormat
$pf.SpaceBefore = 5
$pf.SpaceBeforeAuto = $true
$pf.SpaceAfter = 5
$pf.SpaceAfterAuto = $true
$pf.LineSpacingRule = 0

$nwFont = $normalWeb.Font
$nwFont.NameAscii = "Times New Roman"
$nwFont.NameFarEast = "Times New Roman"
$nwFont.Name = "Times New Roman"
$nwFont.NameBi = "Times New Roman"
$nwFont.Size = 12
$nwFont.SizeBi = 12

# ---------------------------------------------------------------------------
# 2. Append the new paragraphs (blank line, discussion-post answer, sources
#    heading, and the two citation paragraphs, plus a trailing blank line)
#    right after the existing second paragraph, before the final section
#    break. Using InsertXML on a collapsed range at the very end of the
#    document's content lets us specify the exact run/paragraph structure
#    (including the proofErr spell/grammar markers) instead of relying on
#    Word's autoformatting when typing text.
# ---------------------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Hashing is a helpful technique that involves mapping keys and values into a hash table. This allows for faster access to the components of the table. It is also used for other purposes, such as securing passwords to prevent data breaches. Hashing is used in password storage to transform passwords into data that can’t be converted back into the original password. This type of storage prevents attackers from accessing entire tables of usernames and passwords. Instead, they can only access a table of hashed passwords that cannot be converted </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>to</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> original ones. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Hash tables have become a massive part of data scientists' tech interviews because they allow quicker data access and modification. Instead of searching through data for an element, if you know the associated hash key, you can perform searches and changes in constant time, no matter how large the dataset is. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Sources: </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:ind w:left="567" w:hanging="567"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
    <w:t xml:space="preserve">Arias, D. (2019, September 30). </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t xml:space="preserve">How to hash passwords: </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t>One-way road</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t xml:space="preserve"> to enhanced security</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
    <w:t xml:space="preserve">. Auth0. https://auth0.com/blog/hashing-passwords-one-way-road-to-security/ </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:ind w:left="567" w:hanging="567"/>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
    <w:t>Linkedin</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
    <w:t xml:space="preserve">. (2023, July 6). </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t xml:space="preserve">What are some examples of using hash tables to optimize your </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
      <w:i/>
      <w:iCs/>
    </w:rPr>
    <w:t>code?</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/>
    </w:rPr>
    <w:t xml:space="preserve"> Hash Tables: Examples, Challenges, and Interview Tips. https://www.linkedin.com/advice/0/what-some-examples-using-hash-tables-optimize </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$endRange.InsertXML($xml)
